$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.787.02'
$ws.Range('E2').Value = '  +4.33%  '
$ws.Range('D3').Value = '2.284.37'
$ws.Range('E3').Value = '  +2.93%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '322.30'
$ws.Range('E5').Value = '  +1.60%  '
$ws.Range('D6').Value = '106.83'
$ws.Range('E6').Value = '  +6.82%  '
$ws.Range('D7').Value = '0.596'
$ws.Range('E7').Value = '  +0.86%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '0.576'
$ws.Range('E9').Value = '  +2.70%  '
$ws.Range('D10').Value = '38.99'
$ws.Range('E10').Value = '  +5.08%  '
$ws.Range('D11').Value = '0.0847'
$ws.Range('E11').Value = '  +1.81%  '
$ws.Range('D12').Value = '7.95'
$ws.Range('E12').Value = '  +2.72%  '
$ws.Range('D13').Value = '0.108'
$ws.Range('E13').Value = '  +0.93%  '
$ws.Range('D14').Value = '0.890'
$ws.Range('E14').Value = '  +3.47%  '
$ws.Range('D15').Value = '2.631.22'
$ws.Range('E15').Value = '  +2.71%  '
$ws.Range('D16').Value = '14.70'
$ws.Range('E16').Value = '  +3.40%  '
$ws.Range('D17').Value = '2.286.20'
$ws.Range('E17').Value = '  +2.95%  '
$ws.Range('D18').Value = '44.661.46'
$ws.Range('E18').Value = '  +4.19%  '
$ws.Range('D19').Value = '14.13'
$ws.Range('E19').Value = '  -3.39%  '
$ws.Range('E20').Value = '  +4.90%  '
$ws.Range('D21').Value = '6.60'
$ws.Range('E21').Value = '  +3.02%  '
$ws.Range('D22').Value = '66.75'
$ws.Range('E22').Value = '  +2.14%  '
$ws.Range('D23').Value = '3.23'
$ws.Range('E23').Value = '  +2.65%  '
$ws.Range('D24').Value = '240.93'
$ws.Range('E24').Value = '  +2.19%  '
$ws.Range('D25').Value = '2.23'
$ws.Range('E25').Value = '  +3.72%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('D27').Value = '10.25'
$ws.Range('E27').Value = '  +2.51%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').Value = '38.96'
$ws.Range('E28').Value = '  +13.06%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value = '2.23'
$ws.Range('E29').Value = '  +0.77%  '
$ws.Range('D30').Value = '6.56'
$ws.Range('E30').Value = '  +3.80%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '20.76'
$ws.Range('E31').Value = '  +0.95%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '0.0893'
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').Value = '162.92'
$ws.Range('E33').Value = '  +4.58%  '
$ws.Range('D34').Value = '2.80'
$ws.Range('E34').Value = '  +0.56%  '
$ws.Range('E35').Value = '  +10.82%  '
$ws.Range('D36').Value = '2.05'
$ws.Range('E36').Value = '  +5.28%  '
$ws.Range('E37').Value = '  +0.72%  '
$ws.Range('E38').Value = '  +0.23%  '
$ws.Range('D39').Value = '3.96'
$ws.Range('E39').Value = '  +1.74%  '
$ws.Range('D40').Value = '4.50'
$ws.Range('E40').Value = '  +1.71%  '
$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D41').Value = '15.75'
$ws.Range('E41').Value = '  +25.65%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '0.0331'
$ws.Range('E42').Value = '  +1.64%  '
$ws.Range('E43').Value = '  +0.10%  '
$ws.Range('D44').Value = '1.778.84'
$ws.Range('E44').Value = '  -7.14%  '
$ws.Range('B45').Value = 'BitcoinSV'
$ws.Range('C45').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D45').Value = '88.18'
$ws.Range('E45').Value = '  -0.57%  '
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D46').Value = '0.211'
$ws.Range('E46').Value = '  +1.49%  '
$ws.Range('E47').Value = '  +2.46%  '
$ws.Range('B48').Value = 'ordi'
$ws.Range('C48').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D48').Value = '76.90'
$ws.Range('E48').Value = '  -0.19%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').Value = '61.04'
$ws.Range('E49').Value = '  +0.62%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '1.74'
$ws.Range('E50').Value = '  +8.55%  '
$ws.Range('D51').Value = '105.15'
$ws.Range('E51').Value = '  +2.18%  '
